$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.743.29"
$ws.Range("E2").Value = "  -0.21%  "

$ws.Range("D3").Value = "'2.071.00"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "'232.77"
$ws.Range("E5").Value = "  -1.05%  "

$ws.Range("D6").Value = "'0.623"
$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "'58.23"
$ws.Range("E8").Value = "  -1.55%  "

$ws.Range("D9").Value = "'0.391"
$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("D10").Value = "'0.0784"
$ws.Range("E10").Value = "  -1.14%  "

$ws.Range("D11").Value = "'0.110"
$ws.Range("E11").Value = "  +4.19%  "

$ws.Range("D12").Value = "'2.378.67"
$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("D13").Value = "'14.68"
$ws.Range("E13").Value = "  +0.28%  "

$ws.Range("D14").Value = "'21.03"
$ws.Range("E14").Value = "  -1.89%  "

$ws.Range("D15").Value = "'0.774"
$ws.Range("E15").Value = "  +0.32%  "

$ws.Range("D16").Value = "'5.31"
$ws.Range("E16").Value = "  -0.30%  "

$ws.Range("D17").Value = "'2.071.27"
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").Value = "'37.597.73"
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("D19").Value = "'71.34"
$ws.Range("E19").Value = "  -0.52%  "

$ws.Range("D20").Value = "'6.09"
$ws.Range("E20").Value = "  -2.35%  "

$ws.Range("D21").Value = "'0.0₃0840"
$ws.Range("E21").Value = "  +1.24%  "

$ws.Range("D22").Value = "'228.46"
$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").Value = "'2.38"
$ws.Range("E24").Value = "  -1.29%  "

$ws.Range("D25").Value = "'2.39"
$ws.Range("E25").Value = "  -1.32%  "

$ws.Range("D26").Value = "'9.55"
$ws.Range("E26").Value = "  +5.71%  "

$ws.Range("D27").Value = "'171.09"
$ws.Range("E27").Value = "  +0.15%  "

$ws.Range("D28").Value = "'0.138"
$ws.Range("E28").Value = "  -1.87%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'19.36"
$ws.Range("E29").Value = "  -0.89%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.39"
$ws.Range("E30").Value = "  -3.68%  "

$ws.Range("D31").Value = "'0.121"
$ws.Range("E31").Value = "  +0.84%  "

$ws.Range("D32").Value = "'4.71"
$ws.Range("E32").Value = "  +0.32%  "

$ws.Range("D33").Value = "'0.0630"
$ws.Range("E33").Value = "  -0.25%  "

$ws.Range("D34").Value = "'4.67"
$ws.Range("E34").Value = "  -0.75%  "

$ws.Range("D35").Value = "'2.45"
$ws.Range("E35").Value = "  -3.66%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'1.82"
$ws.Range("E36").Value = "  -0.58%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'3.38"
$ws.Range("E37").Value = "  -4.18%  "

$ws.Range("E38").Value = "  +0.16%  "

$ws.Range("D39").Value = "'5.44"
$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("D40").Value = "'0.0231"
$ws.Range("E40").Value = "  +7.29%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'99.83"
$ws.Range("E41").Value = "  +0.54%  "

$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "'0.0974"
$ws.Range("E42").Value = "  -1.26%  "

$ws.Range("D43").Value = "'2.93"
$ws.Range("E43").Value = "  -0.21%  "

$ws.Range("D44").Value = "'17.17"
$ws.Range("E44").Value = "  +6.95%  "

$ws.Range("D45").Value = "'1.443.62"
$ws.Range("E45").Value = "  -1.70%  "

$ws.Range("D46").Value = "'1.15"
$ws.Range("E46").Value = "  -2.05%  "

$ws.Range("D47").Value = "'1.06"
$ws.Range("E47").Value = "  -0.85%  "

$ws.Range("D48").Value = "'4.09"
$ws.Range("E48").Value = "  -3.34%  "

$ws.Range("D49").Value = "'7.31"
$ws.Range("E49").Value = "  -2.16%  "

$ws.Range("D50").Value = "'2.98"
$ws.Range("E50").Value = "  -2.09%  "

$ws.Range("D51").Value = "'2.264.24"
$ws.Range("E51").Value = "  -0.66%  "
